# Commit: "model and template with unit, description and enum"
#
# The bloodgases sheet is a 3-layer header block:
#   row 1 = field name
#   row 2 = field type (#string / #date / #integer / #float, ...)
#   row 3 = (new) field description / enum, one per column
#
# This edit:
#   1. gives the "Temperature" column (I) a type annotation that also
#      carries the unit ("#integer,  unit:celsius" instead of plain
#      "#integer"),
#   2. adds a brand-new row 3 with French field descriptions for the
#      first 7 columns (A-G); the last 3 columns (H, I, J -
#      WaitingTime/Temperature/Result) get no description, i.e. stay
#      empty for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2: refine the "Temperature" column's type with its unit ----------
$ws.Cells.Item(2, 9).Value  = "#integer,  unit:celsius"   # I2 (Temperature)
$ws.Cells.Item(2, 10).Value = "#float"                     # J2 (Result) - unchanged value

# --- row 3 (new): per-column description / enum ----------------------------
$ws.Cells.Item(3, 1).Value = "#Manipulateur"                    # A3 - Operator
$ws.Cells.Item(3, 2).Value = "#Desc:IdentifiantEchantillon"     # B3 - SampleID
$ws.Cells.Item(3, 3).Value = "#Date"                             # C3 - Date
$ws.Cells.Item(3, 4).Value = "#ModeOderatoireLaboratoire"       # D3 - LaboratoryOperatingMode
$ws.Cells.Item(3, 5).Value = "#AppareilLogicielCritique"        # E3 - CriticalApparatusCriticalSoftware
$ws.Cells.Item(3, 6).Value = "#ProduitCritique"                 # F3 - CriticalProduct
$ws.Cells.Item(3, 7).Value = "#LieuStockageDonneesBrutes"       # G3 - RawDataPathway

# H3, I3, J3 (WaitingTime / Temperature / Result) intentionally carry no
# description text - left blank (nothing to set).
